$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C6").Value = -12.12710000000001
$ws.Range("B7").Value = 5.031599999999999
$ws.Range("A9").Value = -21.87490000000001
$ws.Range("B12").Value = 5.575999999999996
$ws.Range("A13").Value = -22.48559999999999
$ws.Range("B14").Value = 5.732399999999997
$ws.Range("C15").Value = -14.08350000000001
$ws.Range("A16").Value = -21.48749999999997
$ws.Range("A18").Value = -22.28560000000001
$ws.Range("B19").Value = 8.974300000000003
$ws.Range("A20").Value = -19.85399999999999
$ws.Range("A26").Value = -21.20369999999997
$ws.Range("B26").Value = 3.961200000000003
$ws.Range("A27").Value = -21.35269999999997
$ws.Range("B27").Value = 5.2567
$ws.Range("C28").Value = -13.0171
$ws.Range("A29").Value = -21.72130000000001
$ws.Range("B29").Value = 5.534199999999998
$ws.Range("C33").Value = -11.51699999999999
$ws.Range("A35").Value = -19.7831
$ws.Range("C35").Value = -12.6315
$ws.Range("A36").Value = -19.54799999999999
$ws.Range("B37").Value = 8.515400000000003
$ws.Range("B38").Value = 4.597199999999998
$ws.Range("C38").Value = -12.0297
$ws.Range("C43").Value = -14.0745
$ws.Range("C44").Value = -13.50979999999999
$ws.Range("A45").Value = -21.86139999999999
$ws.Range("C45").Value = -13.18189999999999
$ws.Range("B47").Value = 5.5281
$ws.Range("C47").Value = -12.06889999999999
$ws.Range("B51").Value = 5.956300000000001
$ws.Range("C51").Value = -11.82129999999999
$ws.Range("B52").Value = 5.553099999999997
$ws.Range("C54").Value = -13.1545
$ws.Range("A55").Value = -22.4729
$ws.Range("B55").Value = 4.572099999999997
$ws.Range("A57").Value = -22.0042
$ws.Range("C57").Value = -13.28479999999998
$ws.Range("C62").Value = -14.25040000000001
$ws.Range("C63").Value = -11.19950000000001
$ws.Range("C67").Value = -10.91959999999999
$ws.Range("A69").Value = -21.63660000000001
$ws.Range("B69").Value = 5.671999999999993
$ws.Range("B70").Value = 5.833200000000006
$ws.Range("C70").Value = -11.62849999999999
$ws.Range("A76").Value = -22.1645
$ws.Range("B76").Value = 5.144199999999999
$ws.Range("A78").Value = -19.85309999999998
$ws.Range("B81").Value = 5.301699999999999
$ws.Range("C81").Value = -11.75319999999999
$ws.Range("A82").Value = -21.97010000000001
$ws.Range("A83").Value = -21.91569999999999
$ws.Range("B83").Value = 6.622400000000006
$ws.Range("C88").Value = -11.84789999999999
$ws.Range("A93").Value = -20.49439999999998
$ws.Range("B94").Value = 5.236899999999999
$ws.Range("C96").Value = -11.72790000000001
$ws.Range("A97").Value = -21.77060000000001
$ws.Range("C99").Value = -12.1234
$ws.Range("B100").Value = 5.267799999999999
$ws.Range("B102").Value = 8.697300000000007
